# Handover data form: remove the "Reached" column (it's no longer needed)
# and fill in the completed handover values for row 2 (receiver side of the
# handover has now been recorded: condition, approvals, dates and status).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Reached" was column N (14). Delete the whole column so everything to its
# right (ReceiverCondition .. Status) shifts one column left (N..U -> N..T).
$ws.Range("N1").EntireColumn.Delete()

# Row 2 values, column by column, reflecting the completed handover record.
$ws.Range("A2").Value = "2ba1c34d"
$ws.Range("B2").Value = 111111111111
$ws.Range("C2").Value = "Drone Equipment"
$ws.Range("D2").Value = "Name7"
$ws.Range("E2").Value = "Make7"
$ws.Range("F2").Value = "Model7"
$ws.Range("G2").Value = "Serial7"
$ws.Range("H2").Value = "SOI ASSAM"
$ws.Range("I2").Value = "SOI TRIPURA"
$ws.Range("J2").Value = "Umar"
$ws.Range("K2").Value = "Umar"
$ws.Range("L2").Value = "Not OK"
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "Good"
$ws.Range("O2").Value = "-"
$ws.Range("P2").Value = "YES"
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = "2024-06-07 12:54:38"
$ws.Range("S2").Value = "2024-06-07 13:08:01"
$ws.Range("T2").Value = "Rejected"
